$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 378, shifting existing rows 378:400 down to 379:401
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new record
$ws.Cells.Item(378, 1).Value = 10
$ws.Cells.Item(378, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(378, 3).Value = "La Araucanía"
$ws.Cells.Item(378, 4).Value = 44610
$ws.Cells.Item(378, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(378, 5).Value = 9
$ws.Cells.Item(378, 6).Value = "Fruta"
$ws.Cells.Item(378, 7).Value = 100108
$ws.Cells.Item(378, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(378, 9).Value = 100108005
$ws.Cells.Item(378, 10).Value = "Piña"
$ws.Cells.Item(378, 11).Value = "Caramelo"
$ws.Cells.Item(378, 12).Value = "Primera"
$ws.Cells.Item(378, 13).Value = 200
$ws.Cells.Item(378, 14).Value = 19000
$ws.Cells.Item(378, 15).Value = 19000
$ws.Cells.Item(378, 16).Value = 19000
$ws.Cells.Item(378, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(378, 18).Value = "Ecuador"
$ws.Cells.Item(378, 19).Value = 1583
$ws.Cells.Item(378, 20).Value = 12
